# Insert two new weekly observation rows (991 and 992) above the existing
# "Cilantro / Vega Central Mapocho de Santiago" data block, pushing the
# previous rows 991-1108 down to 993-1110.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 991..1108 down by two rows.
$ws.Rows("991:992").Insert()

# --- New row 991 ------------------------------------------------------
$ws.Range("A991").Value = 9
$ws.Range("B991").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C991").Value = "Metropolitana"
$ws.Range("D991").Value = 45212
$ws.Range("E991").Value = 13
$ws.Range("F991").Value = 100112040
$ws.Range("G991").Value = "Cilantro"
$ws.Range("H991").Value = "Sin especificar"
$ws.Range("I991").Value = "Primera"
$ws.Range("J991").Value = 70
$ws.Range("K991").Value = 6000
$ws.Range("L991").Value = 6000
$ws.Range("M991").Value = 6000
$ws.Range("N991").Value = "`$/caja 36 atados"
$ws.Range("O991").Value = "Región Metropolitana"
$ws.Range("P991").Value = 167
$ws.Range("Q991").Value = 36
$ws.Range("R991").Value = "Hortaliza"

# --- New row 992 ------------------------------------------------------
$ws.Range("A992").Value = 9
$ws.Range("B992").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C992").Value = "Metropolitana"
$ws.Range("D992").Value = 45212
$ws.Range("E992").Value = 13
$ws.Range("F992").Value = 100112040
$ws.Range("G992").Value = "Cilantro"
$ws.Range("H992").Value = "Sin especificar"
$ws.Range("I992").Value = "Primera"
$ws.Range("J992").Value = 160
$ws.Range("K992").Value = 9000
$ws.Range("L992").Value = 10000
$ws.Range("M992").Value = 9500
$ws.Range("N992").Value = "`$/docena de atados"
$ws.Range("O992").Value = "Región Metropolitana"
$ws.Range("P992").Value = 3167
$ws.Range("Q992").Value = 3
$ws.Range("R992").Value = "Hortaliza"
